# SSDQMvsISO21001 v1.0; wskaźniki v0.1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ISO20001_4-10")
$ws.Activate()

# --- Content edits on sheet "ISO20001_4-10" ---

# B21: reworded note (first phase qualifier added)
$ws.Range("B21").Value = "całość pierwszej fazy tego dotyczy, szczególnie 6."

# Previously-empty B cells that now carry a value
$ws.Range("B109").Value = "brak"
$ws.Range("B122").Value = "`"9.7"
$ws.Range("B123").Value = "brak"
$ws.Range("B128").Value = "`"9.6"

# A6 loses its (invisible) fill-applied style -> back to plain/default formatting
$ws.Range("A6").ClearFormats()

# --- View/layout tweaks ---
$ws.Columns("B").ColumnWidth = 24.27
$excel.ActiveWindow.Zoom = 120
$ws.Range("A18").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E31").Select()
